$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the PanCK segment/aoi labels.
# Order matters: "slide name" needs to land at shared-string index 266,
# followed by PanCK_neg / PanCK_neg-aoi-001 / PanCK_pos / PanCK_pos-aoi-001,
# which is what we get when "PanCK-" is replaced before "PanCK+".
$ws.Cells.Replace("PanCK-", "PanCK_neg")
$ws.Cells.Replace("PanCK+", "PanCK_pos")

# F51 picks up a distinct (but visually identical) style in the source
# edit - reassert the font color theme to force a new style record.
$ws.Range("F51").Font.ThemeColor = 1

# Widen column F to fit the longer "PanCK_xxx-aoi-001" labels.
$ws.Columns("F").ColumnWidth = 47

# Scroll the view down and move the selection to C50.
$ws.Range("C50").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 39
$win.ScrollColumn = 1
